{"js": "// The document contains a centered date title followed by a 20x5 table of\n// simple arithmetic equations (one \"<w:t>\" run per cell). The target\n// revision replaces every single text run, in document order, with a new\n// value (new date + 100 new equations) while leaving all formatting /\n// structure untouched. Since every old value is unique within the\n// document, a straightforward search-and-replace per pair is safe and\n// order independent.\nconst pairs = [\n  [\"2024-04-09 Tuesday\", \"2024-04-10 Wednesday\"],\n  [\"42+15=\", \"63-44=\"],\n  [\"82-5=\", \"44+9=\"],\n  [\"65-46=\", \"47+15=\"],\n  [\"58-5=\", \"34+64=\"],\n  [\"67+28=\", \"22-0=\"],\n  [\"47+20=\", \"50+11=\"],\n  [\"79-49=\", \"31+22=\"],\n  [\"94-81=\", \"72-33=\"],\n  [\"12+4=\", \"12+11=\"],\n  [\"80-10=\", \"7+63=\"],\n  [\"2+83=\", \"20+78=\"],\n  [\"30+14=\", \"97-79=\"],\n  [\"72-39=\", \"75-68=\"],\n  [\"83-15=\", \"35-33=\"],\n  [\"83-28=\", \"30+40=\"],\n  [\"60-2=\", \"3+26=\"],\n  [\"56+17=\", \"9+61=\"],\n  [\"19-8=\", \"65+8=\"],\n  [\"58-48=\", \"98-19=\"],\n  [\"95-75=\", \"3+29=\"],\n  [\"40+49=\", \"53-46=\"],\n  [\"21+4=\", \"43+42=\"],\n  [\"68+20=\", \"55+22=\"],\n  [\"58-31=\", \"84-75=\"],\n  [\"60-25=\", \"82+6=\"],\n  [\"49+8=\", \"67+8=\"],\n  [\"17+3=\", \"65-38=\"],\n  [\"84-19=\", \"13+85=\"],\n  [\"70+21=\", \"3+56=\"],\n  [\"35+24=\", \"43+33=\"],\n  [\"24+19=\", \"22-9=\"],\n  [\"76-33=\", \"71-37=\"],\n  [\"38+11=\", \"68-18=\"],\n  [\"29+48=\", \"43-29=\"],\n  [\"19+3=\", \"60+22=\"],\n  [\"5-5=\", \"95-91=\"],\n  [\"73+22=\", \"93-2=\"],\n  [\"63+3=\", \"72+7=\"],\n  [\"89-62=\", \"30+47=\"],\n  [\"5+48=\", \"39+6=\"],\n  [\"23+34=\", \"47-26=\"],\n  [\"74-23=\", \"21+76=\"],\n  [\"36-4=\", \"78-73=\"],\n  [\"49-36=\", \"77-30=\"],\n  [\"93-1=\", \"13+64=\"],\n  [\"7+87=\", \"58+6=\"],\n  [\"71-40=\", \"37+56=\"],\n  [\"52-28=\", \"26-22=\"],\n  [\"70-1=\", \"59+24=\"],\n  [\"50+36=\", \"44+7=\"],\n  [\"45+46=\", \"72-53=\"],\n  [\"35+28=\", \"83-2=\"],\n  [\"67+24=\", \"40+20=\"],\n  [\"81+4=\", \"27+21=\"],\n  [\"44+2=\", \"4+78=\"],\n  [\"84-70=\", \"48+26=\"],\n  [\"57-19=\", \"39-30=\"],\n  [\"95-54=\", \"0+11=\"],\n  [\"37+6=\", \"61+1=\"],\n  [\"77-24=\", \"31-12=\"],\n  [\"49-2=\", \"48+45=\"],\n  [\"4+79=\", \"23+10=\"],\n  [\"45+34=\", \"24-0=\"],\n  [\"73-64=\", \"81-77=\"],\n  [\"54-48=\", \"71+6=\"],\n  [\"38+10=\", \"76-57=\"],\n  [\"10+66=\", \"92-92=\"],\n  [\"38-15=\", \"78-17=\"],\n  [\"36-35=\", \"70+5=\"],\n  [\"17+23=\", \"60+21=\"],\n  [\"76-49=\", \"60+17=\"],\n  [\"81-5=\", \"21+44=\"],\n  [\"0+93=\", \"18+0=\"],\n  [\"40-24=\", \"28-24=\"],\n  [\"50-0=\", \"5+5=\"],\n  [\"2+78=\", \"6+24=\"],\n  [\"82-19=\", \"64-29=\"],\n  [\"9+57=\", \"25-11=\"],\n  [\"68-65=\", \"4+64=\"],\n  [\"72-51=\", \"12+18=\"],\n  [\"48-31=\", \"3+68=\"],\n  [\"21+66=\", \"7+7=\"],\n  [\"25-13=\", \"61+6=\"],\n  [\"25+8=\", \"79-16=\"],\n  [\"61-29=\", \"14+29=\"],\n  [\"38-12=\", \"26+22=\"],\n  [\"32+58=\", \"85-72=\"],\n  [\"62+12=\", \"58-43=\"],\n  [\"4+39=\", \"45+5=\"],\n  [\"28+7=\", \"59-9=\"],\n  [\"87-82=\", \"66-47=\"],\n  [\"91-58=\", \"26-15=\"],\n  [\"8+72=\", \"42+16=\"],\n  [\"10+82=\", \"12-3=\"],\n  [\"43-17=\", \"60-42=\"],\n  [\"58-2=\", \"52-4=\"],\n  [\"67+20=\", \"85-63=\"],\n  [\"6+73=\", \"3+74=\"],\n  [\"11+9=\", \"76-43=\"],\n  [\"56-49=\", \"8+58=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains a centered date title followed by a 20x5 table of\n# simple arithmetic equations (one run of text per cell). The target\n# revision replaces every single text run, in document order, with a new\n# value (new date + 100 new equations) while leaving all formatting /\n# structure untouched. Since every \"old\" value is unique within the\n# document, a straightforward Find/Replace per pair (exact match, no\n# wildcards) is safe and gives a 1:1 correspondence with the target diff.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2024-04-09 Tuesday\"\n$find.Replacement.Text = \"2024-04-10 Wednesday\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"42+15=\"\n$find.Replacement.Text = \"63-44=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"82-5=\"\n$find.Replacement.Text = \"44+9=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"65-46=\"\n$find.Replacement.Text = \"47+15=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"58-5=\"\n$find.Replacement.Text = \"34+64=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"67+28=\"\n$find.Replacement.Text = \"22-0=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"47+20=\"\n$find.Replacement.Text = \"50+11=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"79-49=\"\n$find.Replacement.Text = \"31+22=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"94-81=\"\n$find.Replacement.Text = \"72-33=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"12+4=\"\n$find.Replacement.Text = \"12+11=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"80-10=\"\n$find.Replacement.Text = \"7+63=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2+83=\"\n$find.Replacement.Text = \"20+78=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"30+14=\"\n$find.Replacement.Text = \"97-79=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"72-39=\"\n$find.Replacement.Text = \"75-68=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"83-15=\"\n$find.Replacement.Text = \"35-33=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"83-28=\"\n$find.Replacement.Text = \"30+40=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"60-2=\"\n$find.Replacement.Text = \"3+26=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"56+17=\"\n$find.Replacement.Text = \"9+61=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"19-8=\"\n$find.Replacement.Text = \"65+8=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"58-48=\"\n$find.Replacement.Text = \"98-19=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"95-75=\"\n$find.Replacement.Text = \"3+29=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"40+49=\"\n$find.Replacement.Text = \"53-46=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"21+4=\"\n$find.Replacement.Text = \"43+42=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"68+20=\"\n$find.Replacement.Text = \"55+22=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"58-31=\"\n$find.Replacement.Text = \"84-75=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"60-25=\"\n$find.Replacement.Text = \"82+6=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"49+8=\"\n$find.Replacement.Text = \"67+8=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"17+3=\"\n$find.Replacement.Text = \"65-38=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"84-19=\"\n$find.Replacement.Text = \"13+85=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"70+21=\"\n$find.Replacement.Text = \"3+56=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"35+24=\"\n$find.Replacement.Text = \"43+33=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"24+19=\"\n$find.Replacement.Text = \"22-9=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"76-33=\"\n$find.Replacement.Text = \"71-37=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"38+11=\"\n$find.Replacement.Text = \"68-18=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"29+48=\"\n$find.Replacement.Text = \"43-29=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"19+3=\"\n$find.Replacement.Text = \"60+22=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"5-5=\"\n$find.Replacement.Text = \"95-91=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"73+22=\"\n$find.Replacement.Text = \"93-2=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"63+3=\"\n$find.Replacement.Text = \"72+7=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"89-62=\"\n$find.Replacement.Text = \"30+47=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"5+48=\"\n$find.Replacement.Text = \"39+6=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"23+34=\"\n$find.Replacement.Text = \"47-26=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"74-23=\"\n$find.Replacement.Text = \"21+76=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"36-4=\"\n$find.Replacement.Text = \"78-73=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"49-36=\"\n$find.Replacement.Text = \"77-30=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"93-1=\"\n$find.Replacement.Text = \"13+64=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"7+87=\"\n$find.Replacement.Text = \"58+6=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"71-40=\"\n$find.Replacement.Text = \"37+56=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"52-28=\"\n$find.Replacement.Text = \"26-22=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"70-1=\"\n$find.Replacement.Text = \"59+24=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"50+36=\"\n$find.Replacement.Text = \"44+7=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"45+46=\"\n$find.Replacement.Text = \"72-53=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"35+28=\"\n$find.Replacement.Text = \"83-2=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"67+24=\"\n$find.Replacement.Text = \"40+20=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"81+4=\"\n$find.Replacement.Text = \"27+21=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"44+2=\"\n$find.Replacement.Text = \"4+78=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"84-70=\"\n$find.Replacement.Text = \"48+26=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"57-19=\"\n$find.Replacement.Text = \"39-30=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"95-54=\"\n$find.Replacement.Text = \"0+11=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"37+6=\"\n$find.Replacement.Text = \"61+1=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"77-24=\"\n$find.Replacement.Text = \"31-12=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"49-2=\"\n$find.Replacement.Text = \"48+45=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"4+79=\"\n$find.Replacement.Text = \"23+10=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"45+34=\"\n$find.Replacement.Text = \"24-0=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"73-64=\"\n$find.Replacement.Text = \"81-77=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"54-48=\"\n$find.Replacement.Text = \"71+6=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"38+10=\"\n$find.Replacement.Text = \"76-57=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"10+66=\"\n$find.Replacement.Text = \"92-92=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"38-15=\"\n$find.Replacement.Text = \"78-17=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"36-35=\"\n$find.Replacement.Text = \"70+5=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"17+23=\"\n$find.Replacement.Text = \"60+21=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"76-49=\"\n$find.Replacement.Text = \"60+17=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"81-5=\"\n$find.Replacement.Text = \"21+44=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"0+93=\"\n$find.Replacement.Text = \"18+0=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"40-24=\"\n$find.Replacement.Text = \"28-24=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"50-0=\"\n$find.Replacement.Text = \"5+5=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2+78=\"\n$find.Replacement.Text = \"6+24=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"82-19=\"\n$find.Replacement.Text = \"64-29=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"9+57=\"\n$find.Replacement.Text = \"25-11=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"68-65=\"\n$find.Replacement.Text = \"4+64=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"72-51=\"\n$find.Replacement.Text = \"12+18=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"48-31=\"\n$find.Replacement.Text = \"3+68=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"21+66=\"\n$find.Replacement.Text = \"7+7=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"25-13=\"\n$find.Replacement.Text = \"61+6=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"25+8=\"\n$find.Replacement.Text = \"79-16=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"61-29=\"\n$find.Replacement.Text = \"14+29=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"38-12=\"\n$find.Replacement.Text = \"26+22=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"32+58=\"\n$find.Replacement.Text = \"85-72=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"62+12=\"\n$find.Replacement.Text = \"58-43=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"4+39=\"\n$find.Replacement.Text = \"45+5=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"28+7=\"\n$find.Replacement.Text = \"59-9=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"87-82=\"\n$find.Replacement.Text = \"66-47=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"91-58=\"\n$find.Replacement.Text = \"26-15=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"8+72=\"\n$find.Replacement.Text = \"42+16=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"10+82=\"\n$find.Replacement.Text = \"12-3=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"43-17=\"\n$find.Replacement.Text = \"60-42=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"58-2=\"\n$find.Replacement.Text = \"52-4=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"67+20=\"\n$find.Replacement.Text = \"85-63=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"6+73=\"\n$find.Replacement.Text = \"3+74=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"11+9=\"\n$find.Replacement.Text = \"76-43=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"56-49=\"\n$find.Replacement.Text = \"8+58=\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
